$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 2..153

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)  # Column G ("Recorded By")
    $val = $cell.Text
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ', '
        if ($parts[0] -eq "System" -and $parts.Length -gt 1) {
            $rest = $parts[1..($parts.Length - 1)]
            $newParts = $rest + ,"System"
            $newVal = $newParts -join ', '
            $cell.Value = $newVal
        }
    }
}
